$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new column F header "Tolerance" and the per-row tolerance values
$ws.Range("F1").Value = "Tolerance"
$ws.Range("F2").Value = "High"
$ws.Range("F3").Value = "Low"
$ws.Range("F4").Value = "High"
$ws.Range("F5").Value = "High"

# Update selection to match the post-edit state (active cell F5)
$ws.Range("F5").Select()
